# Weekly update: insert a new observation row at row 20 (above the row that
# currently holds the 2022-02-10 record), shifting all subsequent rows down
# by one. The sheet's used range grows from A1:R100 to A1:R101.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 20; existing rows 20-100 move to 21-101.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new market observation.
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44811
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100114007
$ws.Range("G20").Value = "Jengibre"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 11638
$ws.Range("N20").Value = "$/caja 13 kilos"
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 895
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = "Hortaliza"
